$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 81

# Text-like columns: force text storage, then strip the formatting stamp
# so the cell keeps default style (no explicit NumberFormat survives).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-20"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "22:25:45"
$ws.Cells.Item($row, 3).Value = "Saturday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "02"
$ws.Cells.Item($row, 4).ClearFormats()

# Numeric columns
$ws.Cells.Item($row, 5).Value = 138789
$ws.Cells.Item($row, 6).Value = 140874
$ws.Cells.Item($row, 7).Value = 171753
$ws.Cells.Item($row, 8).Value = 148847
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 123017
$ws.Cells.Item($row, 11).Value = 223819
$ws.Cells.Item($row, 12).Value = 255683
$ws.Cells.Item($row, 13).Value = 185393
$ws.Cells.Item($row, 14).Value = 110411
$ws.Cells.Item($row, 15).Value = 41269
$ws.Cells.Item($row, 16).Value = 30923
$ws.Cells.Item($row, 17).Value = 73643
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42998
$ws.Cells.Item($row, 20).Value = -1
